# Weekly update for "Hortaliza, Vega Modelo de Temuco - Arveja Verde"
# Two new daily records are inserted at the top of the data table (rows 42-43,
# just before the existing row that used to be 42), pushing all the older
# rows down by two positions. The sheet's used range grows from A1:R67 to
# A1:R69.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows above the current row 42 - this shifts every row
# that used to be 42..67 down to 44..69, carrying their data/styles with
# them automatically.
$ws.Rows("42:43").Insert()

# New row 42
$ws.Range("A42").Value = 10
$ws.Range("B42").Value = "Vega Modelo de Temuco"
$ws.Range("C42").Value = "La Araucanía"
$ws.Range("D42").Value = 44529
$ws.Range("E42").Value = 9
$ws.Range("F42").Value = 100112022
$ws.Range("G42").Value = "Arveja Verde"
$ws.Range("H42").Value = "Sin especificar"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 110
$ws.Range("K42").Value = 15000
$ws.Range("L42").Value = 15000
$ws.Range("M42").Value = 15000
$ws.Range("N42").Value = "$/saco 25 kilos"
$ws.Range("O42").Value = "Región de La Araucanía"
$ws.Range("P42").Value = 600
$ws.Range("Q42").Value = 25
$ws.Range("R42").Value = "Hortaliza"

# New row 43
$ws.Range("A43").Value = 10
$ws.Range("B43").Value = "Vega Modelo de Temuco"
$ws.Range("C43").Value = "La Araucanía"
$ws.Range("D43").Value = 44529
$ws.Range("E43").Value = 9
$ws.Range("F43").Value = 100112022
$ws.Range("G43").Value = "Arveja Verde"
$ws.Range("H43").Value = "Sin especificar"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 210
$ws.Range("K43").Value = 15000
$ws.Range("L43").Value = 15000
$ws.Range("M43").Value = 15000
$ws.Range("N43").Value = "$/saco 25 kilos"
$ws.Range("O43").Value = "Región del Maule"
$ws.Range("P43").Value = 600
$ws.Range("Q43").Value = 25
$ws.Range("R43").Value = "Hortaliza"
